$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values (plain numbers) replacing the old shared-string entries.
$values = @(73749, 73753, 73763, 73764, 73771, 73801, 73682, 73688, 73698, 73699, 73733)

# Clear the previous used range (A1:A21) below the header first.
$ws.Range("A2:A21").ClearContents()

# Write header + new numeric values into A1:A12.
$ws.Range("A1").Value = "DocEntry"
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Column A width shrinks to fit the shorter numeric content (mimics the
# bestFit/AutoFit recalculation Excel performs once the long text values
# are replaced by short 5-digit numbers).
$ws.Columns.Item(1).ColumnWidth = 8

# Update the active selection to B5.
$ws.Range("B5").Select()
